# "2 search test cases added-B47 and B48"
# Adds TestCase_E19, TestCase_E20 and TestCase_E21 rows to the "Test Cases"
# sheet and flips the Results column (PASS -> SKIP) for the pre-existing
# rows 2-19 (those tests are now skipped while the new watchlist cases are
# exercised).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- 1. Existing rows 2-19: Results column PASS -> SKIP -------------------
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 5).Value = "SKIP"
}

# --- 2. New rows 20-22: copy formatting from an existing data row ---------
# (row 5 carries the border + wrap-text format used by the "Description"
# column for every data row) so the new cells land on the same style
# indexes (A/B/D/E -> plain bordered style, C -> wrap-text bordered style)
# instead of minting new ones.
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A20:E22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(22).RowHeight = 45

# Row 20 - TestCase_E19
$ws.Cells.Item(20, 1).Value = "TestCase_E19"
$ws.Cells.Item(20, 2).Value = "OPQA-288"
$ws.Cells.Item(20, 3).Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Cells.Item(20, 4).Value = "Y"
$ws.Cells.Item(20, 5).Value = "SKIP"

# Row 21 - TestCase_E20
$ws.Cells.Item(21, 1).Value = "TestCase_E20"
$ws.Cells.Item(21, 2).Value = "OPQA-290"
$ws.Cells.Item(21, 3).Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Cells.Item(21, 4).Value = "Y"
$ws.Cells.Item(21, 5).Value = "SKIP"

# Row 22 - TestCase_E21
$ws.Cells.Item(22, 1).Value = "TestCase_E21"
$ws.Cells.Item(22, 2).Value = "OPQA-291"
$ws.Cells.Item(22, 3).Value = "Verify that following fields are getting displayed for each post in the watchlist page:`na)Likes`nb)Comments"
$ws.Cells.Item(22, 4).Value = "Y"
$ws.Cells.Item(22, 5).Value = "PASS"

# --- 3. Selection / view ----------------------------------------------------
$ws.Range("D2:D22").Select() | Out-Null

# --- 4. Window size (cosmetic) ---------------------------------------------
$excel.ActiveWindow.Width = 14310
$excel.ActiveWindow.Height = 4995
